# Generate Report for Handoff
#
# Inserts a new "Ready for handoff" entry for file
# 6550b08e-f945-4eb6-81ea-8aec86d39a59 ahead of the existing
# ddc650a5-6cb0-4195-b437-e4d2e34184a2 entry on every sheet
# (Overview, zh-cn, de-de), pushing the ".localization-config"
# row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop existing hyperlinks (rebuilt fully at the end) before shifting rows.
$ws1.Hyperlinks.Delete()

# Insert a fresh row below the current "ddc650a5..." row so the new
# row inherits that row's cell formatting (hyperlink-style column A).
$ws1.Rows.Item(3).Insert()

# Row 3 becomes the (pre-existing) ddc650a5 entry, row 2 becomes the
# brand new 6550b08e entry.
$ws1.Range("A3").Value = "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A2").Value = "6550b08e-f945-4eb6-81ea-8aec86d39a59.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Row 4 (previously row 3) already holds ".localization-config" /
# "Not to be localized" / "Not to be localized" and kept its style.

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/6550b08e-f945-4eb6-81ea-8aec86d39a59.md", "", "", "6550b08e-f945-4eb6-81ea-8aec86d39a59.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/ddc650a5-6cb0-4195-b437-e4d2e34184a2.md", "", "", "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Insert()

# Row 3: existing ddc650a5 entry (values unchanged, just shifted down).
$ws2.Range("A3").Value = "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-09 08:07:20"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

# Row 2: new 6550b08e entry.
$ws2.Range("A2").Value = "6550b08e-f945-4eb6-81ea-8aec86d39a59.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-09 08:07:56"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

# Row 4 (previously row 3) already holds the ".localization-config" entry.

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/6550b08e-f945-4eb6-81ea-8aec86d39a59.md", "", "", "6550b08e-f945-4eb6-81ea-8aec86d39a59.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3da1c36f67cd3b60b347d584886aaa112f6ffa4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.zh-cn.xlf", "", "", "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/ddc650a5-6cb0-4195-b437-e4d2e34184a2.md", "", "", "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71e79cfa9651c053b79a7759482a8e5d9199537d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.zh-cn.xlf", "", "", "ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Insert()

# Row 3: existing ddc650a5 entry (values unchanged, just shifted down).
$ws3.Range("A3").Value = "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-09 08:07:24"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

# Row 2: new 6550b08e entry.
$ws3.Range("A2").Value = "6550b08e-f945-4eb6-81ea-8aec86d39a59.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-09 08:08:00"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

# Row 4 (previously row 3) already holds the ".localization-config" entry.

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/6550b08e-f945-4eb6-81ea-8aec86d39a59.md", "", "", "6550b08e-f945-4eb6-81ea-8aec86d39a59.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3da1c36f67cd3b60b347d584886aaa112f6ffa4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.de-de.xlf", "", "", "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/e2e/ddc650a5-6cb0-4195-b437-e4d2e34184a2.md", "", "", "ddc650a5-6cb0-4195-b437-e4d2e34184a2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/265eb4ca7ea11064b144c151707e5afcbeb0ffc5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.de-de.xlf", "", "", "ddc650a5-6cb0-4195-b437-e4d2e34184a2.1c6863cc9e07637706048818994cc8068a792f25.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/58913db89e5bdec1a628f81573141aa7471e2ca4/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report updated: inserted 6550b08e-f945-4eb6-81ea-8aec86d39a59 handoff rows on Overview, zh-cn, de-de."
